$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.533.99"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.344.44"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'190.45"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").Value = "'564.92"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "3.335.25"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'0.187"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'0.590"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "'48.05"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "'8.73"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "3.874.02"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'612.63"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "'18.19"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "66.537.34"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "3.354.55"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'11.22"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "'18.57"
$ws.Range("E23").Value = "  +10.84%  "
$ws.Range("D24").Value = "'5.18"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "'101.73"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("D28").Value = "'9.79"
$ws.Range("E28").Value = "  +5.77%  "
$ws.Range("D29").Value = "'8.73"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'30.59"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("D32").Value = "'4.10"
$ws.Range("E32").Value = "  +9.57%  "
$ws.Range("D33").Value = "'571.61"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "3.759.26"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").Value = "'57.54"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.0₃0734"
$ws.Range("E39").Value = "  +4.08%  "
$ws.Range("D40").Value = "'34.23"
$ws.Range("E40").Value = "  +7.27%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.33"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.132"
$ws.Range("E42").Value = "  +4.89%  "
$ws.Range("E43").Value = "  +4.84%  "
$ws.Range("D44").Value = "'3.43"
$ws.Range("E44").Value = "  +7.59%  "
$ws.Range("D45").Value = "'0.345"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("D47").Value = "'3.22"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +3.49%  "
